$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 1 (Title): consolidate "A" + " " + "slide" into a single run "A slide"
$tr1 = $s.Shapes.Item(1).TextFrame.TextRange
$tr1.Delete()
$tr1.Text = "A slide"

# Shape 4 (TextBox): consolidate "Followed" + " " + "by" + " " + "a" + " " + "picture"
# into a single run "Followed by a picture"
$tr4 = $s.Shapes.Item(4).TextFrame.TextRange
$tr4.Delete()
$tr4.Text = "Followed by a picture"
